$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "other" attributes of row 2 and row 4 (D, M, N, O, P, R, S)
$row2 = @{
    D = $ws.Range("D2").Value()
    M = $ws.Range("M2").Value()
    N = $ws.Range("N2").Value()
    O = $ws.Range("O2").Value()
    P = $ws.Range("P2").Value()
    R = $ws.Range("R2").Value()
    S = $ws.Range("S2").Value()
}

$row4 = @{
    D = $ws.Range("D4").Value()
    M = $ws.Range("M4").Value()
    N = $ws.Range("N4").Value()
    O = $ws.Range("O4").Value()
    P = $ws.Range("P4").Value()
    R = $ws.Range("R4").Value()
    S = $ws.Range("S4").Value()
}

$ws.Range("D2").Value = $row4.D
$ws.Range("M2").Value = $row4.M
$ws.Range("N2").Value = $row4.N
$ws.Range("O2").Value = $row4.O
$ws.Range("P2").Value = $row4.P
$ws.Range("R2").Value = $row4.R
$ws.Range("S2").Value = $row4.S

$ws.Range("D4").Value = $row2.D
$ws.Range("M4").Value = $row2.M
$ws.Range("N4").Value = $row2.N
$ws.Range("O4").Value = $row2.O
$ws.Range("P4").Value = $row2.P
$ws.Range("R4").Value = $row2.R
$ws.Range("S4").Value = $row2.S

$wb.Save()
